# Applies the edits described by the commit:
#  - refresh the cached "datetimeFigureOut" date field text on the slide
#    master, every slide layout, and the notes master (1/15/17 -> 1/16/17)
#  - reword a few runs of text on slide 2 (AR / MA bullets)

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "1/15/17") {
                $tr.Text = "1/16/17"
            }
        }
    }
}

# Slide master date placeholder
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's date placeholder
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master date placeholder
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes

# Slide 2: reword the AR / MA bullet text
$slide2 = $p.Slides.Item(2)
$body = $slide2.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

$full = $tr.Text
$i1 = $full.IndexOf("AR part ")
if ($i1 -ge 0) {
    $tr.Characters($i1 + 1, 8).Text = "Auto-regression (AR) "
}

$full = $tr.Text
$i2 = $full.IndexOf("MA -")
if ($i2 -ge 0) {
    $tr.Characters($i2 + 1, 4).Text = "Moving Average (MA) -"
}

$full = $tr.Text
$oldShock = "Shocks are measured by moving average of the model."
$i3 = $full.IndexOf($oldShock)
if ($i3 -ge 0) {
    $tr.Characters($i3 + 1, $oldShock.Length).Text = "Market shocks are measured by moving average of the model."
}

Write-Output "Edits applied."
